$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Last status check" timestamp in F1
$ws.Range("F1").Value = "Last status check on: 04.02.2022 15:45"

# Update row 2 values
$ws.Range("B2").Value = 34.9
$ws.Range("C2").Value = 34.5

# D2 becomes a text cell holding "+0.4" (not a number) with default style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "+0.4"
$ws.Range("D2").Style = "Normal"

# E2 becomes a text cell holding the timestamp string (not a date-number) with default style
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2022-02-04 15:45:08"
$ws.Range("E2").Style = "Normal"
